# BOT; UPDATE DATA
# Adds one new daily-report row to each of the "all", "kobe" and "other"
# sheets (date serial 43973 = 2020-05-26) and re-points the active tab
# from "kobe" back to "all" (matching the saved selection state in the
# authored workbook).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": insert new row 45 (old row 45 footer shifts to row 46)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(45).Insert() | Out-Null
$wsAll.Range("A45").Value = 43973
$wsAll.Range("B45").Value = 285
$wsAll.Range("C45").Value = 282
$wsAll.Range("D45").Value = 37
$wsAll.Range("E45").Value = 32
$wsAll.Range("F45").Value = 5
$wsAll.Range("G45").Value = 11
$wsAll.Range("H45").Value = 234

# ---------------------------------------------------------------------
# Sheet "kobe": insert new row 100 (old row 100 footer shifts to row 101)
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows.Item(100).Insert() | Out-Null
$wsKobe.Range("A100").Value = 43973
$wsKobe.Range("B100").Value = 0
$wsKobe.Range("C100").Value = 2982
$wsKobe.Range("D100").Value = 0
$wsKobe.Range("E100").Value = 285
$wsKobe.Range("F100").Value = 32
$wsKobe.Range("G100").Value = 28
$wsKobe.Range("H100").Value = 4
$wsKobe.Range("I100").Value = 11
$wsKobe.Range("J100").Value = 225

# ---------------------------------------------------------------------
# Sheet "other": insert new row 75 (old row 75 footer shifts to row 76)
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(75).Insert() | Out-Null
$wsOther.Range("A75").Value = 43973
$wsOther.Range("B75").Value = 0
$wsOther.Range("C75").Value = 14
$wsOther.Range("D75").Value = 5
$wsOther.Range("E75").Value = 4
$wsOther.Range("F75").Value = 1
$wsOther.Range("G75").Value = 0
$wsOther.Range("H75").Value = 9

# ---------------------------------------------------------------------
# Restore each sheet's last selection, then re-activate "all" last so it
# becomes the workbook's active tab again (previously "kobe" was active).
# ---------------------------------------------------------------------
$wsOther.Activate() | Out-Null
$wsOther.Range("B76").Select() | Out-Null

$wsKobe.Activate() | Out-Null
$wsKobe.Range("B101").Select() | Out-Null

$wsAll.Activate() | Out-Null
$wsAll.Range("J45").Select() | Out-Null
